$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename header strings: "<Field>_old" -> "<Field>_FV2310", "<Field>_new" -> "<Field>_FV2404"
$fields = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt $fields.Count; $i++) {
    $colOld = $i + 1          # A..J  -> 1..10
    $colNew = $i + 12         # L..U -> 12..21
    $ws.Cells.Item(1, $colOld).Value = "$($fields[$i])_FV2310"
    $ws.Cells.Item(1, $colNew).Value = "$($fields[$i])_FV2404"
}
# column K (11) keeps its value "diff" - unchanged

# --- 2) Turn the used range into an Excel table (ListObject) with the default style
$rng = $ws.Range("A1:U67")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"

# --- 3) Freeze the header row (split below row 1)
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
